# Updates the cryptocurrency Price (column D) and Volume(1h) (column E) values
# to the latest scraped figures, preserving the original text formatting
# (percent strings padded with spaces, and "dotted" price strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = <new price text, or $null if unchanged>; E = <new volume text, or $null if unchanged> }
$updates = @{
    2 = @{ D = "27.559.46"; E = "  -2.65%  " }
    3 = @{ D = "1.752.04"; E = "  -3.40%  " }
    4 = @{ D = "1.002"; E = "  +0.15%  " }
    5 = @{ D = "323.87"; E = "  -0.48%  " }
    6 = @{ D = "1.001"; E = "  +0.21%  " }
    7 = @{ D = "0.4446"; E = "  +1.63%  " }
    8 = @{ D = "0.3609"; E = "  -1.52%  " }
    9 = @{ D = "0.07483"; E = "  -2.39%  " }
    10 = @{ D = "42.23"; E = "  -5.57%  " }
    11 = @{ D = "1.105"; E = "  -2.96%  " }
    12 = @{ D = "1.001"; E = "  +0.16%  " }
    13 = @{ D = "20.67"; E = "  -6.01%  " }
    14 = @{ D = "6.031"; E = "  -4.34%  " }
    15 = @{ D = "7.179"; E = "  -4.12%  " }
    16 = @{ D = "1.757.07"; E = "  -3.93%  " }
    17 = @{ D = "92.94"; E = "  -2.35%  " }
    18 = @{ D = "0.00001063"; E = "  -1.49%  " }
    19 = @{ D = "0.06418"; E = "  -1.15%  " }
    20 = @{ D = "1.001"; E = "  +0.25%  " }
    21 = @{ D = $null; E = "  -2.05%  " }
    22 = @{ D = "5.838"; E = $null }
    23 = @{ D = "27.602.62"; E = "  -2.54%  " }
    24 = @{ D = "11.25"; E = "  -2.55%  " }
    25 = @{ D = "2.101"; E = "  -0.34%  " }
    26 = @{ D = "162.92"; E = "  +0.93%  " }
    27 = @{ D = "20.40"; E = "  -1.57%  " }
    28 = @{ D = "1.955.80"; E = "  -3.53%  " }
    29 = @{ D = "2.125"; E = "  -6.76%  " }
    30 = @{ D = "125.57"; E = "  -2.64%  " }
    31 = @{ D = "1.084"; E = "  -9.98%  " }
    32 = @{ D = "0.09021"; E = "  -1.37%  " }
    33 = @{ D = "3.637"; E = "  +2.44%  " }
    34 = @{ D = "5.541"; E = "  -7.96%  " }
    35 = @{ D = "12.09"; E = "  -6.56%  " }
    36 = @{ D = "0.02301"; E = "  -2.55%  " }
    37 = @{ D = "0.2094"; E = "  -3.55%  " }
    38 = @{ D = "0.6351"; E = "  -3.57%  " }
    39 = @{ D = "0.05953"; E = "  -4.11%  " }
    40 = @{ D = "4.934"; E = "  -5.40%  " }
    41 = @{ D = "1.193"; E = "  +0.29%  " }
    42 = @{ D = "1.001"; E = "  +0.31%  " }
    43 = @{ D = "1.390"; E = "  -2.50%  " }
    44 = @{ D = "7.793"; E = "  -3.74%  " }
    45 = @{ D = "13.20"; E = "  -4.83%  " }
    46 = @{ D = "3.714"; E = "  -0.52%  " }
    47 = @{ D = "0.5867"; E = "  -4.03%  " }
    48 = @{ D = $null; E = "  -2.79%  " }
    49 = @{ D = "121.53"; E = "  -3.08%  " }
    50 = @{ D = "1.157"; E = "  +0.29%  " }
    51 = @{ D = "0.06850"; E = "  -2.07%  " }
}

foreach ($row in $updates.Keys) {
    $u = $updates[$row]
    if ($null -ne $u.D) {
        # Prefix with an apostrophe so Excel stores the new price as literal text
        # instead of re-interpreting it as a number (which would corrupt values like
        # "27.559.46" or strip significant trailing zeros from values like "13.20").
        $ws.Cells.Item($row, 4).Value = "'" + $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
